# The "31 -35" label in A2 is missing a space before "35" compared to the
# other range labels in the column (e.g. "36 - 40"). Fix the text and
# leave the selection on the corrected cell, matching the target state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "31 - 35"

$ws.Range("A2").Select()
